$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(121, 8).Value = 849.52  # H121: 5028.769 -> 849.52
$ws.Cells.Item(121, 10).Value = 849.52  # J121: 5028.769 -> 849.52
$ws.Cells.Item(121, 12).Value = 2548.56  # L121: 15086.307 -> 2548.56
$ws.Cells.Item(121, 14).Value = -6042.559999999999  # N121: -18580.307 -> -6042.559999999999
$ws.Cells.Item(129, 8).Value = 846.2857  # H129: 847.38336 -> 846.2857
$ws.Cells.Item(129, 9).Value = 597.2  # I129: 597.6 -> 597.2
$ws.Cells.Item(129, 10).Value = 870.7059  # J129: 870.0909 -> 870.7059
$ws.Cells.Item(129, 11).Value = 1791.6  # K129: 1792.8 -> 1791.6
$ws.Cells.Item(129, 12).Value = 2612.1177  # L129: 2610.2727 -> 2612.1177
$ws.Cells.Item(129, 13).Value = 3208.4  # M129: 3207.2 -> 3208.4
$ws.Cells.Item(129, 14).Value = -12612.1177  # N129: -12610.2727 -> -12612.1177
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 28028.342  # H32: 25524.578 -> 28028.342
$ws.Cells.Item(32, 9).Value = 34666.098  # I32: 30688.371 -> 34666.098
$ws.Cells.Item(32, 11).Value = 34666.098  # K32: 30688.371 -> 34666.098
$ws.Cells.Item(32, 13).Value = -34379.098  # M32: -30401.371 -> -34379.098
$ws.Cells.Item(74, 8).Value = 712.8  # H74: 793.5806 -> 712.8
$ws.Cells.Item(74, 9).Value = 361.46155  # I74: 406.86365 -> 361.46155
$ws.Cells.Item(74, 10).Value = 1727.7778  # J74: 1738.8889 -> 1727.7778
$ws.Cells.Item(74, 11).Value = 361.46155  # K74: 406.86365 -> 361.46155
$ws.Cells.Item(74, 12).Value = 1727.7778  # L74: 1738.8889 -> 1727.7778
$ws.Cells.Item(74, 13).Value = 512.53845  # M74: 467.13635 -> 512.53845
$ws.Cells.Item(74, 14).Value = -3475.7778  # N74: -3486.8889 -> -3475.7778
$ws.Cells.Item(77, 8).Value = 712.8  # H77: 793.5806 -> 712.8
$ws.Cells.Item(77, 9).Value = 361.46155  # I77: 406.86365 -> 361.46155
$ws.Cells.Item(77, 10).Value = 1727.7778  # J77: 1738.8889 -> 1727.7778
$ws.Cells.Item(77, 11).Value = 1807.30775  # K77: 2034.31825 -> 1807.30775
$ws.Cells.Item(77, 12).Value = 8638.889000000001  # L77: 8694.4445 -> 8638.889000000001
$ws.Cells.Item(77, 13).Value = 2560.69225  # M77: 2333.68175 -> 2560.69225
$ws.Cells.Item(77, 14).Value = -17374.889  # N77: -17430.4445 -> -17374.889
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 771.1739  # H80: 952.1539 -> 771.1739
$ws.Cells.Item(80, 9).Value = 1606  # I80: 1706 -> 1606
$ws.Cells.Item(80, 10).Value = 476.52942  # J80: 617.1111 -> 476.52942
$ws.Cells.Item(80, 11).Value = 1606  # K80: 1706 -> 1606
$ws.Cells.Item(80, 12).Value = 476.52942  # L80: 617.1111 -> 476.52942
$ws.Cells.Item(80, 13).Value = -608  # M80: -708 -> -608
$ws.Cells.Item(80, 14).Value = -2472.52942  # N80: -2613.1111 -> -2472.52942
$ws.Cells.Item(83, 8).Value = 771.1739  # H83: 952.1539 -> 771.1739
$ws.Cells.Item(83, 9).Value = 1606  # I83: 1706 -> 1606
$ws.Cells.Item(83, 10).Value = 476.52942  # J83: 617.1111 -> 476.52942
$ws.Cells.Item(83, 11).Value = 8030  # K83: 8530 -> 8030
$ws.Cells.Item(83, 12).Value = 2382.6471  # L83: 3085.5555 -> 2382.6471
$ws.Cells.Item(83, 13).Value = -3038  # M83: -3538 -> -3038
$ws.Cells.Item(83, 14).Value = -12366.6471  # N83: -13069.5555 -> -12366.6471
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1180  # H16: 1139.9166 -> 1180
$ws.Cells.Item(16, 9).Value = 1100  # I16: 1068.4286 -> 1100
$ws.Cells.Item(16, 10).Value = 1233.3334  # J16: 1240 -> 1233.3334
$ws.Cells.Item(16, 11).Value = 1100  # K16: 1068.4286 -> 1100
$ws.Cells.Item(16, 12).Value = 1233.3334  # L16: 1240 -> 1233.3334
$ws.Cells.Item(16, 13).Value = -813  # M16: -781.4286 -> -813
$ws.Cells.Item(16, 14).Value = -1807.3334  # N16: -1814 -> -1807.3334
$ws.Cells.Item(80, 8).Value = 23000  # H80: 21000 -> 23000
$ws.Cells.Item(80, 10).Value = 23000  # J80: 21000 -> 23000
$ws.Cells.Item(80, 12).Value = 23000  # L80: 21000 -> 23000
$ws.Cells.Item(80, 14).Value = -25246  # N80: -23246 -> -25246
$ws.Cells.Item(83, 8).Value = 23000  # H83: 21000 -> 23000
$ws.Cells.Item(83, 10).Value = 23000  # J83: 21000 -> 23000
$ws.Cells.Item(83, 12).Value = 69000  # L83: 63000 -> 69000
$ws.Cells.Item(83, 14).Value = -80232  # N83: -74232 -> -80232
$ws.Cells.Item(99, 8).Value = 5656.1113  # H99: 5505.263 -> 5656.1113
$ws.Cells.Item(99, 9).Value = 4081  # I99: 3960 -> 4081
$ws.Cells.Item(99, 10).Value = 7625  # J99: 7222.222 -> 7625
$ws.Cells.Item(99, 11).Value = 4081  # K99: 3960 -> 4081
$ws.Cells.Item(99, 12).Value = 7625  # L99: 7222.222 -> 7625
$ws.Cells.Item(99, 13).Value = -2583  # M99: -2462 -> -2583
$ws.Cells.Item(99, 14).Value = -10621  # N99: -10218.222 -> -10621
$ws.Cells.Item(113, 8).Value = 1180  # H113: 1139.9166 -> 1180
$ws.Cells.Item(113, 9).Value = 1100  # I113: 1068.4286 -> 1100
$ws.Cells.Item(113, 10).Value = 1233.3334  # J113: 1240 -> 1233.3334
$ws.Cells.Item(113, 11).Value = 1100  # K113: 1068.4286 -> 1100
$ws.Cells.Item(113, 12).Value = 1233.3334  # L113: 1240 -> 1233.3334
$ws.Cells.Item(113, 13).Value = 1070  # M113: 1101.5714 -> 1070
$ws.Cells.Item(113, 14).Value = -5573.3334  # N113: -5580 -> -5573.3334
$ws.Cells.Item(122, 8).Value = 990  # H122: 766.1429000000001 -> 990
$ws.Cells.Item(122, 9).Value = 1399.8572  # I122: 895.2727 -> 1399.8572
$ws.Cells.Item(122, 10).Value = 272.75  # J122: 292.66666 -> 272.75
$ws.Cells.Item(122, 11).Value = 4199.571599999999  # K122: 2685.8181 -> 4199.571599999999
$ws.Cells.Item(122, 12).Value = 818.25  # L122: 877.9999799999999 -> 818.25
$ws.Cells.Item(122, 13).Value = -1749.571599999999  # M122: -235.8181 -> -1749.571599999999
$ws.Cells.Item(122, 14).Value = -5718.25  # N122: -5777.99998 -> -5718.25
$ws.Cells.Item(126, 8).Value = 5656.1113  # H126: 5505.263 -> 5656.1113
$ws.Cells.Item(126, 9).Value = 4081  # I126: 3960 -> 4081
$ws.Cells.Item(126, 10).Value = 7625  # J126: 7222.222 -> 7625
$ws.Cells.Item(126, 11).Value = 12243  # K126: 11880 -> 12243
$ws.Cells.Item(126, 12).Value = 22875  # L126: 21666.666 -> 22875
$ws.Cells.Item(126, 13).Value = -9773  # M126: -9410 -> -9773
$ws.Cells.Item(126, 14).Value = -27815  # N126: -26606.666 -> -27815
$ws.Cells.Item(134, 8).Value = 533.4167  # H134: 663.6875 -> 533.4167
$ws.Cells.Item(134, 9).Value = 523.8570999999999  # I134: 657.9286 -> 523.8570999999999
$ws.Cells.Item(134, 10).Value = 600.3333  # J134: 704 -> 600.3333
$ws.Cells.Item(134, 11).Value = 1571.5713  # K134: 1973.7858 -> 1571.5713
$ws.Cells.Item(134, 12).Value = 1800.9999  # L134: 2112 -> 1800.9999
$ws.Cells.Item(134, 13).Value = 963.4287000000002  # M134: 561.2142000000001 -> 963.4287000000002
$ws.Cells.Item(134, 14).Value = -6870.9999  # N134: -7182 -> -6870.9999
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 653.9729599999999  # H5: 675.7838 -> 653.9729599999999
$ws.Cells.Item(5, 9).Value = 477.05264  # I5: 523.58826 -> 477.05264
$ws.Cells.Item(5, 10).Value = 840.7222  # J5: 805.15 -> 840.7222
$ws.Cells.Item(5, 11).Value = 1431.15792  # K5: 1570.76478 -> 1431.15792
$ws.Cells.Item(5, 12).Value = 2522.1666  # L5: 2415.45 -> 2522.1666
$ws.Cells.Item(5, 13).Value = -1319.15792  # M5: -1458.76478 -> -1319.15792
$ws.Cells.Item(5, 14).Value = -2746.1666  # N5: -2639.45 -> -2746.1666
$ws.Cells.Item(18, 8).Value = 186  # H18: 206.66667 -> 186
$ws.Cells.Item(18, 9).Value = 120  # I18: 158 -> 120
$ws.Cells.Item(18, 11).Value = 360  # K18: 474 -> 360
$ws.Cells.Item(18, 13).Value = -191  # M18: -305 -> -191
$ws.Cells.Item(23, 8).Value = 587.4  # H23: 697.1111 -> 587.4
$ws.Cells.Item(23, 9).Value = 600.6667  # I23: 500.5 -> 600.6667
$ws.Cells.Item(23, 10).Value = 585.05884  # J23: 753.2857 -> 585.05884
$ws.Cells.Item(23, 11).Value = 1802.0001  # K23: 1501.5 -> 1802.0001
$ws.Cells.Item(23, 12).Value = 1755.17652  # L23: 2259.8571 -> 1755.17652
$ws.Cells.Item(23, 13).Value = -1567.0001  # M23: -1266.5 -> -1567.0001
$ws.Cells.Item(23, 14).Value = -2225.17652  # N23: -2729.8571 -> -2225.17652
$ws.Cells.Item(131, 8).Value = 189528.8  # H131: 164781.03 -> 189528.8
$ws.Cells.Item(131, 10).Value = 204917.67  # J131: 176273.39 -> 204917.67
$ws.Cells.Item(131, 12).Value = 614753.01  # L131: 528820.17 -> 614753.01
$ws.Cells.Item(131, 14).Value = -624833.01  # N131: -538900.17 -> -624833.01
$ws.Cells.Item(135, 8).Value = 653.9729599999999  # H135: 675.7838 -> 653.9729599999999
$ws.Cells.Item(135, 9).Value = 477.05264  # I135: 523.58826 -> 477.05264
$ws.Cells.Item(135, 10).Value = 840.7222  # J135: 805.15 -> 840.7222
$ws.Cells.Item(135, 11).Value = 4293.47376  # K135: 4712.29434 -> 4293.47376
$ws.Cells.Item(135, 12).Value = 7566.499800000001  # L135: 7246.349999999999 -> 7566.499800000001
$ws.Cells.Item(135, 13).Value = -1758.47376  # M135: -2177.29434 -> -1758.47376
$ws.Cells.Item(135, 14).Value = -12636.4998  # N135: -12316.35 -> -12636.4998
$ws.Cells.Item(139, 8).Value = 11769.777  # H139: 14942.571 -> 11769.777
$ws.Cells.Item(139, 9).Value = 13116  # I139: 14942.571 -> 13116
$ws.Cells.Item(139, 10).Value = 1000  # J139: 0 -> 1000
$ws.Cells.Item(139, 11).Value = 39348  # K139: 44827.713 -> 39348
$ws.Cells.Item(139, 12).Value = 3000  # L139: 0 -> 3000
$ws.Cells.Item(139, 13).Value = -34208  # M139: -39687.713 -> -34208
$ws.Cells.Item(139, 14).Value = -13280  # N139: None -> -13280
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4304.1816  # H70: 4354.7 -> 4304.1816
$ws.Cells.Item(70, 9).Value = 3932.8333  # I70: 3959.6 -> 3932.8333
$ws.Cells.Item(70, 11).Value = 3932.8333  # K70: 3959.6 -> 3932.8333
$ws.Cells.Item(70, 13).Value = -3662.8333  # M70: -3689.6 -> -3662.8333
$ws.Cells.Item(73, 8).Value = 4304.1816  # H73: 4354.7 -> 4304.1816
$ws.Cells.Item(73, 9).Value = 3932.8333  # I73: 3959.6 -> 3932.8333
$ws.Cells.Item(73, 11).Value = 3932.8333  # K73: 3959.6 -> 3932.8333
$ws.Cells.Item(73, 13).Value = -2996.8333  # M73: -3023.6 -> -2996.8333
$ws.Cells.Item(113, 8).Value = 3428.524  # H113: 3231.5833 -> 3428.524
$ws.Cells.Item(113, 9).Value = 2776.4119  # I113: 2708.7778 -> 2776.4119
$ws.Cells.Item(113, 10).Value = 6200  # J113: 4800 -> 6200
$ws.Cells.Item(113, 11).Value = 2776.4119  # K113: 2708.7778 -> 2776.4119
$ws.Cells.Item(113, 12).Value = 6200  # L113: 4800 -> 6200
$ws.Cells.Item(113, 13).Value = -606.4119000000001  # M113: -538.7777999999998 -> -606.4119000000001
$ws.Cells.Item(113, 14).Value = -10540  # N113: -9140 -> -10540
$ws.Cells.Item(122, 8).Value = 3443.7727  # H122: 3573.15 -> 3443.7727
$ws.Cells.Item(122, 9).Value = 2417.0667  # I122: 2482.5715 -> 2417.0667
$ws.Cells.Item(122, 10).Value = 5643.857  # J122: 6117.8335 -> 5643.857
$ws.Cells.Item(122, 11).Value = 7251.2001  # K122: 7447.7145 -> 7251.2001
$ws.Cells.Item(122, 12).Value = 16931.571  # L122: 18353.5005 -> 16931.571
$ws.Cells.Item(122, 13).Value = -4801.2001  # M122: -4997.7145 -> -4801.2001
$ws.Cells.Item(122, 14).Value = -21831.571  # N122: -23253.5005 -> -21831.571
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 952.13635  # H16: 957.8182 -> 952.13635
$ws.Cells.Item(16, 9).Value = 725.2778  # I16: 732.2222 -> 725.2778
$ws.Cells.Item(16, 11).Value = 725.2778  # K16: 732.2222 -> 725.2778
$ws.Cells.Item(16, 13).Value = -555.2778  # M16: -562.2222 -> -555.2778
$ws.Cells.Item(40, 8).Value = 4884.857  # H40: 3816.1667 -> 4884.857
$ws.Cells.Item(40, 9).Value = 4238.8  # I40: 3279.4 -> 4238.8
$ws.Cells.Item(40, 11).Value = 4238.8  # K40: 3279.4 -> 4238.8
$ws.Cells.Item(40, 13).Value = -4102.8  # M40: -3143.4 -> -4102.8
$ws.Cells.Item(64, 8).Value = 0  # H64: 9800 -> 0
$ws.Cells.Item(64, 10).Value = 0  # J64: 9800 -> 0
$ws.Cells.Item(64, 12).Value = 0  # L64: 9800 -> 0
$ws.Cells.Item(64, 14).Value = $null  # N64: -10250 -> None
$ws.Cells.Item(67, 8).Value = 0  # H67: 9800 -> 0
$ws.Cells.Item(67, 10).Value = 0  # J67: 9800 -> 0
$ws.Cells.Item(67, 12).Value = 0  # L67: 9800 -> 0
$ws.Cells.Item(67, 14).Value = $null  # N67: -11360 -> None
$ws.Cells.Item(122, 8).Value = 3225.6428  # H122: 3355.7273 -> 3225.6428
$ws.Cells.Item(122, 9).Value = 2736.2856  # I122: 2727 -> 2736.2856
$ws.Cells.Item(122, 11).Value = 8208.856800000001  # K122: 8181 -> 8208.856800000001
$ws.Cells.Item(122, 13).Value = -5758.856800000001  # M122: -5731 -> -5758.856800000001
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1338.6666  # H107: 1597.421 -> 1338.6666
$ws.Cells.Item(107, 9).Value = 681.3333  # I107: 737.75 -> 681.3333
$ws.Cells.Item(107, 10).Value = 1831.6666  # J107: 2222.6365 -> 1831.6666
$ws.Cells.Item(107, 11).Value = 2043.9999  # K107: 2213.25 -> 2043.9999
$ws.Cells.Item(107, 12).Value = 5494.9998  # L107: 6667.9095 -> 5494.9998
$ws.Cells.Item(107, 13).Value = -123.9999  # M107: -293.25 -> -123.9999
$ws.Cells.Item(107, 14).Value = -9334.9998  # N107: -10507.9095 -> -9334.9998
$ws.Cells.Item(122, 8).Value = 2034.5625  # H122: 2054.5 -> 2034.5625
$ws.Cells.Item(122, 9).Value = 1815.4  # I122: 1794.8889 -> 1815.4
$ws.Cells.Item(122, 10).Value = 2399.8333  # J122: 2833.3333 -> 2399.8333
$ws.Cells.Item(122, 11).Value = 5446.200000000001  # K122: 5384.6667 -> 5446.200000000001
$ws.Cells.Item(122, 12).Value = 7199.499899999999  # L122: 8499.999899999999 -> 7199.499899999999
$ws.Cells.Item(122, 13).Value = -2996.200000000001  # M122: -2934.6667 -> -2996.200000000001
$ws.Cells.Item(122, 14).Value = -12099.4999  # N122: -13399.9999 -> -12099.4999
$ws.Cells.Item(131, 8).Value = 26999  # H131: 0 -> 26999
$ws.Cells.Item(131, 10).Value = 26999  # J131: 0 -> 26999
$ws.Cells.Item(131, 12).Value = 26999  # L131: 0 -> 26999
$ws.Cells.Item(131, 14).Value = -37079  # N131: None -> -37079
$ws.Cells.Item(133, 8).Value = 0  # H133: 1715 -> 0
$ws.Cells.Item(133, 10).Value = 0  # J133: 1715 -> 0
$ws.Cells.Item(133, 12).Value = 0  # L133: 1715 -> 0
$ws.Cells.Item(133, 14).Value = $null  # N133: -11835 -> None
